$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.649.17"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "1.891.93"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'238.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4835"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("D8").Value = "'0.2864"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'0.06554"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "1.972.67"
$ws.Range("E10").Value = "  +6.34%  "
$ws.Range("D11").Value = "'0.07473"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").Value = "'16.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("D13").Value = "'5.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "'0.6672"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("D16").Value = "30.624.80"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "2.215.48"
$ws.Range("E19").Value = "  +5.67%  "
$ws.Range("D20").Value = "'0.000007578"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "'230.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'6.205"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("D25").Value = "'169.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("D26").Value = "'9.364"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("D27").Value = "'18.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "'1.962"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").Value = "'0.1025"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.72%  "
$ws.Range("D30").Value = "'1.400"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("D32").Value = "'4.026"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").Value = "'0.05058"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "'1.215"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.24%  "
$ws.Range("D35").Value = "'0.7548"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.27%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'2.712"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "'0.01872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").Value = "'0.9213"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").Value = "'2.066"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").Value = "'107.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "'0.4291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "'5.669"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("D46").Value = "'7.418"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "'64.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "'0.1275"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "'8.987"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").Value = "'33.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
